$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "E" column cross-referencing the transliterated names of the
# first half of the Coptic year (Tut..Nasi) alongside the existing Coptic
# month names in column C / D.
$ws.Range("E3").Value  = "Tut"
$ws.Range("E4").Value  = "Babah"
$ws.Range("E5").Value  = "Hatur"
$ws.Range("E6").Value  = "Kiahk"
$ws.Range("E7").Value  = "Tubah"
$ws.Range("E8").Value  = "Amshir"
$ws.Range("E9").Value  = "Baramhat"
$ws.Range("E10").Value = "Baramouda"
$ws.Range("E11").Value = "Bashons"
$ws.Range("E12").Value = "Baunah"
$ws.Range("E13").Value = "Abib"
$ws.Range("E14").Value = "Mesra"
$ws.Range("E15").Value = "Nasi"

# Move the active cell selection to E9, matching the saved view state.
$ws.Range("E9").Select()

# Row 3 reverts to the sheet's default (auto-fit) height in the saved file.
$ws.Rows.Item(3).AutoFit()

# Page orientation is explicitly recorded as portrait in the saved file.
$ws.PageSetup.Orientation = 1
